$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before J (TANGGAL PERAWATAN and everything to its right
# shifts one column to the right) to make room for the new "KETERANGAN" column.
$null = $ws.Columns("J:J").Insert()
$ws.Columns("J:J").ColumnWidth = 21.25

# Rename the "SATUAN" header to "STATUS KEPEMILIKAN" (column G) - the
# underlying row-2 formula for that column (${table:pk.SATUAN}) is unchanged.
$ws.Range("G1").Value = "STATUS KEPEMILIKAN"

# Populate the newly inserted column J with its header + data-binding formula.
$ws.Range("J1").Value = "KETERANGAN"
$ws.Range("J2").Value = "`${table:pk.KETERANGAN}"

# Give the new column's body cell the same look as its neighbours (centered,
# bordered) plus word-wrap so longer remarks are readable.
$ws.Range("J2").HorizontalAlignment = -4108
$ws.Range("J2").WrapText = $true

# Cosmetic row-height tweak that came along with the edit.
$ws.Rows("2:2").RowHeight = 17.5

# Restore the selection to the newly added cell, matching the saved view state.
$null = $ws.Range("J2").Select()
